# Rename ObjTables metadata attributes from PascalCase to lowerCamelCase
# (ObjTablesVersion -> objTablesVersion, Type -> type, Id -> id) in the
# special header strings stored in cell A1 (and A2 on the table-of-contents
# sheet) of each metadata worksheet.

$wb = $excel.ActiveWorkbook

$toc = $wb.Worksheets.Item("!!_Table of contents")
$toc.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$toc.Range("A2").Value = "!!ObjTables type='TableOfContents'"

$dataRepo = $wb.Worksheets.Item("!!Data repo metadata")
$dataRepo.Range("A1").Value = "!!ObjTables type='Data' id='DataRepoMetadata'"

$schemaRepo = $wb.Worksheets.Item("!!Schema repo metadata")
$schemaRepo.Range("A1").Value = "!!ObjTables type='Data' id='SchemaRepoMetadata'"

$model1s = $wb.Worksheets.Item("!!Model1s")
$model1s.Range("A1").Value = "!!ObjTables type='Data' id='Model1'"
